# Applies the numeric-value corrections described by the commit diff.
# All target cells are plain numeric literals (no formulas in this workbook).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1914.1666
$ws.Range("I2").Value = 829
$ws.Range("J2").Value = 2999.3333
$ws.Range("K2").Value = 829
$ws.Range("L2").Value = 2999.3333
$ws.Range("M2").Value = -716
$ws.Range("N2").Value = -3225.3333
$ws.Range("H4").Value = 1081.9286
$ws.Range("I4").Value = 249.57143
$ws.Range("J4").Value = 1914.2858
$ws.Range("K4").Value = 249.57143
$ws.Range("L4").Value = 1914.2858
$ws.Range("M4").Value = -135.57143
$ws.Range("N4").Value = -2142.2858
$ws.Range("H5").Value = 97.40000000000001
$ws.Range("I5").Value = 97.40000000000001
$ws.Range("K5").Value = 97.40000000000001
$ws.Range("M5").Value = 17.59999999999999
$ws.Range("H6").Value = 501.54544
$ws.Range("I6").Value = 451.7
$ws.Range("K6").Value = 1355.1
$ws.Range("M6").Value = -1243.1
$ws.Range("H17").Value = 31590.637
$ws.Range("J17").Value = 33565.547
$ws.Range("L17").Value = 100696.641
$ws.Range("N17").Value = -101032.641
$ws.Range("H41").Value = 931.73334
$ws.Range("I41").Value = 1361.75
$ws.Range("J41").Value = 440.2857
$ws.Range("K41").Value = 1361.75
$ws.Range("L41").Value = 440.2857
$ws.Range("M41").Value = -921.75
$ws.Range("N41").Value = -1320.2857
$ws.Range("H70").Value = 1222977
$ws.Range("J70").Value = 3516.5
$ws.Range("L70").Value = 10549.5
$ws.Range("N70").Value = -11089.5
$ws.Range("H73").Value = 1222977
$ws.Range("J73").Value = 3516.5
$ws.Range("L73").Value = 10549.5
$ws.Range("N73").Value = -12421.5
$ws.Range("H76").Value = 7100.846
$ws.Range("I76").Value = 8113.1
$ws.Range("K76").Value = 8113.1
$ws.Range("M76").Value = -7798.1
$ws.Range("H79").Value = 7100.846
$ws.Range("I79").Value = 8113.1
$ws.Range("K79").Value = 8113.1
$ws.Range("M79").Value = -7021.1
$ws.Range("H107").Value = 1144.1111
$ws.Range("I107").Value = 379.5
$ws.Range("J107").Value = 3820.25
$ws.Range("K107").Value = 379.5
$ws.Range("L107").Value = 3820.25
$ws.Range("M107").Value = 1540.5
$ws.Range("N107").Value = -7660.25
$ws.Range("H137").Value = 16668988
$ws.Range("J137").Value = 2687.7058
$ws.Range("L137").Value = 8063.117400000001
$ws.Range("N137").Value = -13163.1174
$ws.Range("H138").Value = 3699.65
$ws.Range("I138").Value = 2436.524
$ws.Range("J138").Value = 4379.795
$ws.Range("K138").Value = 7309.572
$ws.Range("L138").Value = 13139.385
$ws.Range("M138").Value = -2169.572
$ws.Range("N138").Value = -23419.385

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2751.4167
$ws.Range("J2").Value = 2892.2856
$ws.Range("L2").Value = 2892.2856
$ws.Range("N2").Value = -3118.2856
$ws.Range("H32").Value = 4080.0625
$ws.Range("I32").Value = 4080.0625
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 4080.0625
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -3793.0625
$ws.Range("N32").Value = $null
$ws.Range("H45").Value = 5216.4
$ws.Range("J45").Value = 7039.5557
$ws.Range("L45").Value = 7039.5557
$ws.Range("N45").Value = -7793.5557
$ws.Range("H60").Value = 88848.09
$ws.Range("I60").Value = 88848.09
$ws.Range("K60").Value = 88848.09
$ws.Range("M60").Value = -88115.09
$ws.Range("H61").Value = 14494761
$ws.Range("I61").Value = 17077726
$ws.Range("K61").Value = 17077726
$ws.Range("M61").Value = -17077514
$ws.Range("H97").Value = 1151.8518
$ws.Range("I97").Value = 1255.45
$ws.Range("J97").Value = 855.8570999999999
$ws.Range("K97").Value = 1255.45
$ws.Range("L97").Value = 855.8570999999999
$ws.Range("M97").Value = -759.45
$ws.Range("N97").Value = -1847.8571
$ws.Range("H116").Value = 2751.4167
$ws.Range("J116").Value = 2892.2856
$ws.Range("L116").Value = 2892.2856
$ws.Range("N116").Value = -7480.2856
$ws.Range("H132").Value = 1728357.8
$ws.Range("I132").Value = 3462.9363
$ws.Range("J132").Value = 9098363
$ws.Range("K132").Value = 10388.8089
$ws.Range("L132").Value = 27295089
$ws.Range("M132").Value = -7858.8089
$ws.Range("N132").Value = -27300149
$ws.Range("H136").Value = 14494761
$ws.Range("I136").Value = 17077726
$ws.Range("K136").Value = 51233178
$ws.Range("M136").Value = -51230628

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2751.4167
$ws.Range("J3").Value = 2892.2856
$ws.Range("L3").Value = 2892.2856
$ws.Range("N3").Value = -3120.2856
$ws.Range("H22").Value = 662.61536
$ws.Range("I22").Value = 697
$ws.Range("K22").Value = 697
$ws.Range("M22").Value = -524

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").Value = $null
$ws.Range("H31").Value = 26319676
$ws.Range("I31").Value = 41669790
$ws.Range("J31").Value = 5191.2144
$ws.Range("K31").Value = 41669790
$ws.Range("L31").Value = 5191.2144
$ws.Range("M31").Value = -41669495
$ws.Range("N31").Value = -5781.2144
$ws.Range("H34").Value = 26319676
$ws.Range("I34").Value = 41669790
$ws.Range("J34").Value = 5191.2144
$ws.Range("K34").Value = 41669790
$ws.Range("L34").Value = 5191.2144
$ws.Range("M34").Value = -41669588
$ws.Range("N34").Value = -5595.2144
$ws.Range("H62").Value = 15878778
$ws.Range("I62").Value = 4906.533
$ws.Range("K62").Value = 4906.533
$ws.Range("M62").Value = -4282.533
$ws.Range("H65").Value = 15878778
$ws.Range("I65").Value = 4906.533
$ws.Range("K65").Value = 24532.665
$ws.Range("M65").Value = -21412.665
$ws.Range("H125").Value = 87221
$ws.Range("J125").Value = 87221
$ws.Range("L125").Value = 87221
$ws.Range("N125").Value = -92141

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 286.25
$ws.Range("I6").Value = 286.25
$ws.Range("K6").Value = 858.75
$ws.Range("M6").Value = -745.75
$ws.Range("H12").Value = 4562.5
$ws.Range("I12").Value = 13.666667
$ws.Range("J12").Value = 6078.778
$ws.Range("K12").Value = 41.000001
$ws.Range("L12").Value = 18236.334
$ws.Range("M12").Value = 131.999999
$ws.Range("N12").Value = -18582.334
$ws.Range("H107").Value = 4549739.5
$ws.Range("I107").Value = 1917.625
$ws.Range("K107").Value = 5752.875
$ws.Range("M107").Value = -3832.875
$ws.Range("H137").Value = 7824.737
$ws.Range("J137").Value = 11137.7
$ws.Range("L137").Value = 33413.10000000001
$ws.Range("N137").Value = -43613.10000000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 37503230
$ws.Range("J80").Value = 66670176
$ws.Range("L80").Value = 66670176
$ws.Range("N80").Value = -66672172
$ws.Range("H83").Value = 37503230
$ws.Range("J83").Value = 66670176
$ws.Range("L83").Value = 333350880
$ws.Range("N83").Value = -333360864
$ws.Range("H122").Value = 4569.2354
$ws.Range("I122").Value = 3080.182
$ws.Range("J122").Value = 7299.1665
$ws.Range("K122").Value = 9240.545999999998
$ws.Range("L122").Value = 21897.4995
$ws.Range("M122").Value = -6790.545999999998
$ws.Range("N122").Value = -26797.4995
$ws.Range("H132").Value = 6253429
$ws.Range("I132").Value = 3471.2222
$ws.Range("K132").Value = 10413.6666
$ws.Range("M132").Value = -7883.6666

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 722.62164
$ws.Range("I16").Value = 688.13336
$ws.Range("J16").Value = 870.4286
$ws.Range("K16").Value = 688.13336
$ws.Range("L16").Value = 870.4286
$ws.Range("M16").Value = -518.13336
$ws.Range("N16").Value = -1210.4286
$ws.Range("H55").Value = 1430.5294
$ws.Range("I55").Value = 1320
$ws.Range("J55").Value = 1528.7778
$ws.Range("K55").Value = 1320
$ws.Range("L55").Value = 1528.7778
$ws.Range("M55").Value = -1147
$ws.Range("N55").Value = -1874.7778
$ws.Range("H132").Value = 3150.2222
$ws.Range("I132").Value = 2234.7646
$ws.Range("K132").Value = 6704.293799999999
$ws.Range("M132").Value = -4174.293799999999
$ws.Range("H136").Value = 5099.48
$ws.Range("I136").Value = 5180.4443
$ws.Range("K136").Value = 15541.3329
$ws.Range("M136").Value = -12991.3329
